$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

# --- Update existing rows 2-13: new dates (Generic 1 Bonds futures roll dates) and signal values ---

# Row 2
$ws.Range("B2").Value = 43202
$ws.Range("C2").Value = 5.7402213506523632
$ws.Range("D2").Value = -0.97341583846487134
$ws.Range("E2").Value = 1.9044519681442327
$ws.Range("F2").Value = 5.8118930175807035
$ws.Range("G2").Value = 7.5034108265953714

# Row 3
$ws.Range("B3").Value = 43202
$ws.Range("C3").Value = 9.1926037297326726
$ws.Range("D3").Value = 11.071906557705075
$ws.Range("E3").Value = -3.1249848846569099
$ws.Range("F3").Value = 1.5561839012784087
$ws.Range("G3").Value = 7.1870144627716224

# Row 4
$ws.Range("B4").Value = 43202
$ws.Range("C4").Value = 11.893659167996583
$ws.Range("D4").Value = 14.962751586945219
$ws.Range("E4").Value = -3.9864128256911044
$ws.Range("F4").Value = 0.71743931090858926
$ws.Range("G4").Value = 6.0273566265938756

# Row 5
$ws.Range("B5").Value = 43202
$ws.Range("C5").Value = 7.8894100441314787
$ws.Range("D5").Value = 9.0601177926637089
$ws.Range("E5").Value = 6.5107556788116838
$ws.Range("F5").Value = -1.648741041082459
$ws.Range("G5").Value = -1.9900944311644155

# Row 6
$ws.Range("B6").Value = 43202
$ws.Range("C6").Value = 9.4318026430809532
$ws.Range("D6").Value = 8.6643598196580474
$ws.Range("E6").Value = -2.0895562836341606
$ws.Range("F6").Value = 4.7499633479461512
$ws.Range("G6").Value = 11.477478137787582

# Row 7
$ws.Range("B7").Value = 43203
$ws.Range("C7").Value = 0.44848005394504492
$ws.Range("D7").Value = 7.9361794159837489
$ws.Range("E7").Value = 3.4445833664773353
$ws.Range("F7").Value = -12.449710886312944
$ws.Range("G7").Value = -15.318537170169225

# Row 8
$ws.Range("B8").Value = 43203
$ws.Range("C8").Value = 6.2471313095442653
$ws.Range("D8").Value = 10.761804814096337
$ws.Range("E8").Value = 5.5490149701918474
$ws.Range("F8").Value = -7.0063056639735963
$ws.Range("G8").Value = -12.117767770460732

# Row 9
$ws.Range("B9").Value = 43203
$ws.Range("C9").Value = 4.7637205081046554
$ws.Range("D9").Value = 9.5274410162093108
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

# Row 10
$ws.Range("B10").Value = 43203
$ws.Range("C10").Value = 15.690889177423504
$ws.Range("D10").Value = 20
$ws.Range("E10").Value = 10.617588522339645
$ws.Range("F10").Value = -4.5368558884230303
$ws.Range("G10").Value = -8.4037961491734166

# Row 11
$ws.Range("B11").Value = 43202
$ws.Range("C11").Value = 8.6685046173644515
$ws.Range("D11").Value = 3.2817711759426516
$ws.Range("E11").Value = 8.2752801270979859
$ws.Range("F11").Value = 9.4882856806693336
$ws.Range("G11").Value = 12.266451755101885

# Row 12
$ws.Range("B12").Value = 43203
$ws.Range("C12").Value = -1.5573721835490968
$ws.Range("D12").Value = -3.1147443670981936
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0

# Row 13
$ws.Range("B13").Value = 43203
$ws.Range("C13").Value = -16.559879116594324
$ws.Range("D13").Value = -15.532888311392476
$ws.Range("E13").Value = 4.5942347219017865
$ws.Range("F13").Value = -15.811544175405892
$ws.Range("G13").Value = -20

# --- New row 14: DAX ---
$ws.Range("A14").Value = "DAX"
$ws.Range("B14").Value = 43202
$ws.Range("B14").Style = $ws.Range("B13").Style
$ws.Range("C14").Value = -3.0558073289275525
$ws.Range("D14").Value = 0.14393696283560189
$ws.Range("E14").Value = -4.5883094093854924
$ws.Range("F14").Value = -6.6442603023127358
$ws.Range("G14").Value = -0.80845670779221257

# --- Column widths: re-fit after the new (shorter) row 14 label/values ---
$ws.Columns.Item(1).ColumnWidth = 8.7109375
$ws.Columns.Item(3).ColumnWidth = 16.42578125
